# Append the new daily data row (row 77) to Sheet1, matching the
# auto-generated daily push entry for 2025-10-08 01:59 UTC.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds literal text dates like "2025/09/22", not real date
# values, so force text formatting before assigning to avoid Excel
# auto-converting the string into a date serial number; then clear the
# format back off so the cell keeps the workbook's default (unstyled)
# look, matching the other data rows.
$ws.Range("A77").NumberFormat = "@"
$ws.Range("A77").Value = "2025/10/08"
$ws.Range("A77").Style = "Normal"

$ws.Range("B77").Value = "水"
$ws.Range("C77").Value = 10
$ws.Range("D77").Value = 123
